$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 45179 = 2023-09-10) for every
# data row from row 2 through row 375. The edit bumps that date forward by one
# day (serial 45180 = 2023-09-11) across all of those rows.
$ws.Range("C2:C375").Value = 45180
